$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MySQL")

# Bring over the exact formatting (styles) used by the last existing data
# row so the two new rows render identically (wrap text columns B:G,
# date format in H).
$ws.Range("A43:H43").Copy()
$ws.Range("A44:H45").PasteSpecial(-4122)

# --- Row 44: new "guangzhouyihong" base-import table -----------------
$ws.Range("A44").Value = "guangzhouyihong"
$ws.Range("B44").Value = "import"
$ws.Range("C44").Value = "基础"
$ws.Range("D44").Value = "广州一弘基础供货数据"
$ws.Range("E44").Value = "无"
$ws.Range("F44").Value = "无"
$ws.Range("G44").Value = "无"
$ws.Range("H44").Value = 42950

# --- Row 45: new "guangzhouyihong_sichuan_claim_cmp" analysis table --
$ws.Range("A45").Value = "guangzhouyihong_sichuan_claim_cmp"
$ws.Range("B45").Value = "analysis"
$ws.Range("C45").Value = "分析"
$ws.Range("D45").Value = "广州一弘数据与四川理赔数据对比分析"
$ws.Range("E45").Value = "guangzhouyihong_cmp"
$ws.Range("F45").Value = "guangzhouyihong"
$ws.Range("G45").Value = "base_lclaim"
$ws.Range("H45").Value = 42950

# Row 45 wraps onto several lines, same as the other "analysis" rows.
$ws.Rows.Item(45).RowHeight = 40.5

# Scroll the view down to the new rows and leave the selection where the
# author left it after typing in the new data.
$excel.ActiveWindow.ScrollRow = 31
$ws.Range("H46").Select()
